# Update market-price / profit figures in the Leve profit sheets.
# Source data refreshed by the scheduled market-board runner; only the
# computed price/profit columns (H-N) change, one worksheet row at a time.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 80
$ws.Range("H80").Value = 1955.1786
$ws.Range("J80").Value = 2470.5557
$ws.Range("L80").Value = 7411.6671
$ws.Range("N80").Value = -9407.667099999999

# Row 83
$ws.Range("H83").Value = 1955.1786
$ws.Range("J83").Value = 2470.5557
$ws.Range("L83").Value = 22235.0013
$ws.Range("N83").Value = -32219.0013

# Row 101
$ws.Range("H101").Value = 3808.1667
$ws.Range("J101").Value = 4224.75
$ws.Range("L101").Value = 12674.25
$ws.Range("N101").Value = -15918.25

# Row 113
$ws.Range("H113").Value = 6862.972
$ws.Range("I113").Value = 6190.25
$ws.Range("J113").Value = 7703.875
$ws.Range("K113").Value = 6190.25
$ws.Range("L113").Value = 7703.875
$ws.Range("M113").Value = -2936.25
$ws.Range("N113").Value = -14211.875

# Row 131
$ws.Range("H131").Value = 6643.75
$ws.Range("I131").Value = 4138.4614
$ws.Range("J131").Value = 17500
$ws.Range("K131").Value = 12415.3842
$ws.Range("L131").Value = 52500
$ws.Range("M131").Value = -7375.3842
$ws.Range("N131").Value = -62580

# Row 132
$ws.Range("H132").Value = 1827.1364
$ws.Range("I132").Value = 1552.6111
$ws.Range("J132").Value = 3062.5
$ws.Range("K132").Value = 4657.8333
$ws.Range("L132").Value = 9187.5
$ws.Range("M132").Value = -2127.8333
$ws.Range("N132").Value = -14247.5

# Row 135
$ws.Range("H135").Value = 9438160
$ws.Range("I135").Value = 14707485
$ws.Range("K135").Value = 132367365
$ws.Range("M135").Value = -132364830

# Row 137
$ws.Range("H137").Value = 3380.6875
$ws.Range("I137").Value = 3211.5
$ws.Range("K137").Value = 9634.5
$ws.Range("M137").Value = -7084.5

# Row 138
$ws.Range("H138").Value = 6400.0835
$ws.Range("I138").Value = 3561.375
$ws.Range("J138").Value = 7819.4375
$ws.Range("K138").Value = 10684.125
$ws.Range("L138").Value = 23458.3125
$ws.Range("M138").Value = -5544.125
$ws.Range("N138").Value = -33738.3125

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 23576456
$ws.Range("I32").Value = 25662402
$ws.Range("K32").Value = 25662402
$ws.Range("M32").Value = -25662115

# Row 45
$ws.Range("H45").Value = 3756.7
$ws.Range("I45").Value = 1895.875
$ws.Range("K45").Value = 1895.875
$ws.Range("M45").Value = -1518.875

# Row 61
$ws.Range("H61").Value = 5236.8076
$ws.Range("I61").Value = 3150.6667
$ws.Range("J61").Value = 8081.5454
$ws.Range("K61").Value = 3150.6667
$ws.Range("L61").Value = 8081.5454
$ws.Range("M61").Value = -2938.6667
$ws.Range("N61").Value = -8505.545399999999

# Row 74
$ws.Range("H74").Value = 4343.1665
$ws.Range("I74").Value = 4006.75
$ws.Range("K74").Value = 4006.75
$ws.Range("M74").Value = -3132.75

# Row 77
$ws.Range("H77").Value = 4343.1665
$ws.Range("I77").Value = 4006.75
$ws.Range("K77").Value = 20033.75
$ws.Range("M77").Value = -15665.75

# Row 132
$ws.Range("H132").Value = 4337.616
$ws.Range("I132").Value = 3702.4138
$ws.Range("J132").Value = 5653.393
$ws.Range("K132").Value = 11107.2414
$ws.Range("L132").Value = 16960.179
$ws.Range("M132").Value = -8577.241399999999
$ws.Range("N132").Value = -22020.179

# Row 136
$ws.Range("H136").Value = 5236.8076
$ws.Range("I136").Value = 3150.6667
$ws.Range("J136").Value = 8081.5454
$ws.Range("K136").Value = 9452.000100000001
$ws.Range("L136").Value = 24244.6362
$ws.Range("M136").Value = -6902.000100000001
$ws.Range("N136").Value = -29344.6362

$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 11767715
$ws.Range("I94").Value = 2295.5715
$ws.Range("K94").Value = 2295.5715
$ws.Range("M94").Value = -1844.5715

# Row 134
$ws.Range("H134").Value = 4378.972
$ws.Range("I134").Value = 3505.84
$ws.Range("K134").Value = 10517.52
$ws.Range("M134").Value = -7982.52

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 4797.75
$ws.Range("I31").Value = 3411.2083
$ws.Range("J31").Value = 6184.2915
$ws.Range("K31").Value = 3411.2083
$ws.Range("L31").Value = 6184.2915
$ws.Range("M31").Value = -3116.2083
$ws.Range("N31").Value = -6774.2915

# Row 34
$ws.Range("H34").Value = 4797.75
$ws.Range("I34").Value = 3411.2083
$ws.Range("J34").Value = 6184.2915
$ws.Range("K34").Value = 3411.2083
$ws.Range("L34").Value = 6184.2915
$ws.Range("M34").Value = -3209.2083
$ws.Range("N34").Value = -6588.2915

# Row 58
$ws.Range("H58").Value = 3521.1
$ws.Range("I58").Value = 2315.8572
$ws.Range("K58").Value = 2315.8572
$ws.Range("M58").Value = -2112.8572

# Row 70
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()

# Row 73
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()

# Row 86
$ws.Range("H86").Value = 8747.159
$ws.Range("I86").Value = 11888.4
$ws.Range("K86").Value = 11888.4
$ws.Range("M86").Value = -10765.4

# Row 89
$ws.Range("H89").Value = 8747.159
$ws.Range("I89").Value = 11888.4
$ws.Range("K89").Value = 59442
$ws.Range("M89").Value = -53826

# Row 122
$ws.Range("H122").Value = 2947.111
$ws.Range("I122").Value = 1879
$ws.Range("J122").Value = 5083.3335
$ws.Range("K122").Value = 5637
$ws.Range("L122").Value = 15250.0005
$ws.Range("M122").Value = -3187
$ws.Range("N122").Value = -20150.0005

# Row 132
$ws.Range("H132").Value = 4979.129
$ws.Range("I132").Value = 4471.6665
$ws.Range("J132").Value = 5454.875
$ws.Range("K132").Value = 13414.9995
$ws.Range("L132").Value = 16364.625
$ws.Range("M132").Value = -10884.9995
$ws.Range("N132").Value = -21424.625

# Row 136
$ws.Range("H136").Value = 3521.1
$ws.Range("I136").Value = 2315.8572
$ws.Range("K136").Value = 6947.571599999999
$ws.Range("M136").Value = -4397.571599999999

$ws = $wb.Worksheets.Item("CUL")
# Row 3
$ws.Range("H3").Value = 5272.8335
$ws.Range("I3").Value = 5272.8335
$ws.Range("K3").Value = 15818.5005
$ws.Range("M3").Value = -15706.5005

# Row 113
$ws.Range("H113").Value = 3451.4
$ws.Range("J113").Value = 3436.125
$ws.Range("L113").Value = 10308.375
$ws.Range("N113").Value = -14648.375

# Row 129
$ws.Range("H129").Value = 2333.4
$ws.Range("J129").Value = 2333.4
$ws.Range("L129").Value = 7000.200000000001
$ws.Range("N129").Value = -17000.2

# Row 134
$ws.Range("H134").Value = 7019.696
$ws.Range("I134").Value = 2715.8125
$ws.Range("K134").Value = 8147.4375
$ws.Range("M134").Value = -3077.4375

# Row 139
$ws.Range("H139").Value = 33345526
$ws.Range("J139").Value = 14727.272
$ws.Range("L139").Value = 44181.81600000001
$ws.Range("N139").Value = -54461.81600000001

$ws = $wb.Worksheets.Item("GSM")
# Row 113
$ws.Range("H113").Value = 6719.5
$ws.Range("I113").Value = 4083.3333
$ws.Range("J113").Value = 10673.75
$ws.Range("K113").Value = 4083.3333
$ws.Range("L113").Value = 10673.75
$ws.Range("M113").Value = -1913.3333
$ws.Range("N113").Value = -15013.75

# Row 122
$ws.Range("H122").Value = 4147
$ws.Range("I122").Value = 3868.2856
$ws.Range("J122").Value = 5122.5
$ws.Range("K122").Value = 11604.8568
$ws.Range("L122").Value = 15367.5
$ws.Range("M122").Value = -9154.856800000001
$ws.Range("N122").Value = -20267.5

# Row 126
$ws.Range("H126").Value = 2571.6428
$ws.Range("J126").Value = 3382.5715
$ws.Range("L126").Value = 10147.7145
$ws.Range("N126").Value = -15087.7145

# Row 132
$ws.Range("H132").Value = 9217.795
$ws.Range("I132").Value = 7725.7407
$ws.Range("K132").Value = 23177.2221
$ws.Range("M132").Value = -20647.2221

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 5737
$ws.Range("I7").Value = 2801
$ws.Range("J7").Value = 6838
$ws.Range("K7").Value = 2801
$ws.Range("L7").Value = 6838
$ws.Range("M7").Value = -2689
$ws.Range("N7").Value = -7062

# Row 46
$ws.Range("H46").Value = 7605.952
$ws.Range("I46").Value = 4817.6
$ws.Range("J46").Value = 8477.3125
$ws.Range("K46").Value = 4817.6
$ws.Range("L46").Value = 8477.3125
$ws.Range("M46").Value = -4629.6
$ws.Range("N46").Value = -8853.3125

# Row 126
$ws.Range("H126").Value = 5737
$ws.Range("I126").Value = 2801
$ws.Range("J126").Value = 6838
$ws.Range("K126").Value = 8403
$ws.Range("L126").Value = 20514
$ws.Range("M126").Value = -5933
$ws.Range("N126").Value = -25454

# Row 132
$ws.Range("H132").Value = 4498.129
$ws.Range("I132").Value = 3551.9092
$ws.Range("J132").Value = 6811.1113
$ws.Range("K132").Value = 10655.7276
$ws.Range("L132").Value = 20433.3339
$ws.Range("M132").Value = -8125.7276
$ws.Range("N132").Value = -25493.3339

# Row 136
$ws.Range("H136").Value = 5832.7144
$ws.Range("I136").Value = 3985
$ws.Range("J136").Value = 7865.2
$ws.Range("K136").Value = 11955
$ws.Range("L136").Value = 23595.6
$ws.Range("M136").Value = -9405
$ws.Range("N136").Value = -28695.6

$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 4421.524
$ws.Range("I122").Value = 3943.7693
$ws.Range("K122").Value = 11831.3079
$ws.Range("M122").Value = -9381.3079

# Row 132
$ws.Range("H132").Value = 3491
$ws.Range("I132").Value = 2383.1777
$ws.Range("J132").Value = 6260.5557
$ws.Range("K132").Value = 7149.533100000001
$ws.Range("L132").Value = 18781.6671
$ws.Range("M132").Value = -4619.533100000001
$ws.Range("N132").Value = -23841.6671

# Row 136
$ws.Range("H136").Value = 4208.7075
$ws.Range("I136").Value = 4079.8262
$ws.Range("J136").Value = 4373.3887
$ws.Range("K136").Value = 12239.4786
$ws.Range("L136").Value = 13120.1661
$ws.Range("M136").Value = -9689.4786
$ws.Range("N136").Value = -18220.1661
